$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text so numeric-looking values (e.g. "1.006")
# are preserved exactly as text instead of being parsed into numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '28.475.16'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '1.870.88'
$ws.Range("E3").Value = '  -0.51%  '

$ws.Range("E4").Value = '  -1.67%  '

$ws.Range("D5").Value = '315.32'
$ws.Range("E5").Value = '  -0.91%  '

$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -1.67%  '

$ws.Range("E7").Value = '  -1.18%  '

$ws.Range("D8").Value = '0.3896'
$ws.Range("E8").Value = '  -1.28%  '

$ws.Range("D9").Value = '0.08347'
$ws.Range("E9").Value = '  -0.02%  '

$ws.Range("E10").Value = '  -1.44%  '

$ws.Range("E11").Value = '  -0.95%  '

$ws.Range("D12").Value = '6.210'
$ws.Range("E12").Value = '  -0.87%  '

$ws.Range("D13").Value = '1.871.35'
$ws.Range("E13").Value = '  +0.54%  '

$ws.Range("D14").Value = '20.41'
$ws.Range("E14").Value = '  -0.25%  '

$ws.Range("D15").Value = '7.279'
$ws.Range("E15").Value = '  +0.35%  '

$ws.Range("E16").Value = '  -1.64%  '

$ws.Range("D17").Value = '0.00001103'
$ws.Range("E17").Value = '  -0.74%  '

$ws.Range("D18").Value = '91.08'
$ws.Range("E18").Value = '  -0.54%  '

$ws.Range("D19").Value = '0.06727'
$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").Value = '17.72'

$ws.Range("E21").Value = '  -1.73%  '

$ws.Range("D22").Value = '5.905'
$ws.Range("E22").Value = '  -1.22%  '

$ws.Range("D23").Value = '28.495.30'
$ws.Range("E23").Value = '  -0.32%  '

$ws.Range("D24").Value = '11.12'
$ws.Range("E24").Value = '  -0.79%  '

$ws.Range("D25").Value = '2.216'
$ws.Range("E25").Value = '  -2.46%  '

$ws.Range("D26").Value = '2.087.44'

$ws.Range("D27").Value = '160.79'

$ws.Range("D28").Value = '20.63'
$ws.Range("E28").Value = '  -1.01%  '

$ws.Range("D29").Value = '2.403'
$ws.Range("E29").Value = '  +0.72%  '

$ws.Range("D30").Value = '126.31'
$ws.Range("E30").Value = '  -1.06%  '

$ws.Range("E31").Value = '  -0.44%  '

$ws.Range("D32").Value = '1.037'
$ws.Range("E32").Value = '  +0.13%  '

$ws.Range("D33").Value = '5.739'
$ws.Range("E33").Value = '  -1.96%  '

$ws.Range("D34").Value = '3.610'
$ws.Range("E34").Value = '  -1.39%  '

$ws.Range("D35").Value = '0.02449'
$ws.Range("E35").Value = '  +0.40%  '

$ws.Range("D36").Value = '0.06548'
$ws.Range("E36").Value = '  +0.63%  '

$ws.Range("B37").Value = 'Algorand'
$ws.Range("C37").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D37").Value = '0.2163'
$ws.Range("E37").Value = '  -1.14%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").Value = '8.901'
$ws.Range("E38").Value = '  -2.74%  '

$ws.Range("D39").Value = '5.016'
$ws.Range("E39").Value = '  +0.26%  '

$ws.Range("E40").Value = '  -1.00%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.234'
$ws.Range("E41").Value = '  -1.47%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").Value = '0.6365'
$ws.Range("E42").Value = '  -1.32%  '

$ws.Range("E43").Value = '  -1.37%  '

$ws.Range("E44").Value = '  -1.98%  '

$ws.Range("D45").Value = '0.6001'

$ws.Range("D46").Value = '13.02'
$ws.Range("E46").Value = '  -0.07%  '

$ws.Range("D47").Value = '3.689'
$ws.Range("E47").Value = '  -0.90%  '

$ws.Range("D48").Value = '1.999'
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("D49").Value = '1.215'
$ws.Range("E49").Value = '  +0.04%  '

$ws.Range("D50").Value = '121.84'
$ws.Range("E50").Value = '  -0.34%  '

$ws.Range("D51").Value = '1.131'
$ws.Range("E51").Value = '  -10.00%  '

# Restore the default "Normal" style on the price column so no stray
# number-format styling is left behind on the cells.
$priceRange.Style = "Normal"
